# Updates league bases data: for a handful of duplicate-date fixtures,
# the match rows got shuffled back to their correct order. For each pair
# of rows below, swap the entire row contents in columns B:AC (the "id"
# through "PL_AhUnder" fields), leaving column A (the sequential row
# number) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(362, 363),
    @(371, 372),
    @(377, 378),
    @(444, 445)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
